$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The map grid (A1:P20) was redrawn/moved (see commit message: "Moved Map
# to server"). Rebuild the full 20-row x 16-col tile layout in its new
# state and write it back onto the sheet.
$data = @(
    @(2,2,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(2,2,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,1,1,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,1,1,1,1,3,3,3,3,3,3,3,3,3,3),
    @(3,1,1,1,1,1,1,3,3,3,3,3,3,3,3,3),
    @(3,1,1,1,1,1,1,3,3,3,3,3,3,3,3,3),
    @(3,3,1,1,1,1,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,1,1,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,1,1,1,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,0,0,1,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,0,0,1,2,2,2,3,3,3,3,3,3,3,3,3),
    @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3),
    @(3,3,3,3,3,3,3,3,3,3,3,3,3,3,3,3)
)

for ($r = 0; $r -lt 20; $r++) {
    for ($c = 0; $c -lt 16; $c++) {
        $ws.Cells.Item($r + 1, $c + 1).Value = $data[$r][$c]
    }
}

# Restore the active-cell selection used after the edit.
$ws.Range("A2").Select()
